$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows before the existing row 2 so the current 2017-2019
# data (rows 2-7) shifts down to rows 10-15, making room for the new
# 2013-2016 data.
$insertRange = $ws.Range("A2:E9")
$insertRange.EntireRow.Insert()

# New data for years 2013-2016 (Foreign / Retail Domestic pairs)
$newRows = @(
    @(2013, "Foreign", 40, 3.48, 8.7),
    @(2013, "Retail Domestic", 34, 2.958, 8.7),
    @(2014, "Foreign", 47, 4.841, 10.3),
    @(2014, "Retail Domestic", 34, 3.502, 10.3),
    @(2015, "Foreign", 44, 4.136, 9.4),
    @(2015, "Retail Domestic", 34, 3.196, 9.4),
    @(2016, "Foreign", 45, 4.185, 9.3),
    @(2016, "Retail Domestic", 34, 3.162, 9.3)
)

$r = 2
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$ws.Range("D12").Select()
